$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = -3.062305129253262
$ws.Range("F23").Value = -3.073686251417679
$ws.Range("F25").Value = -3.097599818497839
$ws.Range("F26").Value = -3.097571928340787
$ws.Range("F27").Value = -3.112290222692378
$ws.Range("F28").Value = -3.11827770017562
$ws.Range("F29").Value = -3.12529970434338
$ws.Range("F30").Value = -3.132838727428308
$ws.Range("F31").Value = -3.129861307059424
$ws.Range("F32").Value = -3.138798314754496
$ws.Range("F33").Value = -3.152882305271939
$ws.Range("F34").Value = -3.15593283341838
$ws.Range("F35").Value = -3.154106735800806
$ws.Range("F37").Value = -3.163015364488649
$ws.Range("F39").Value = -3.175519919247344
$ws.Range("F40").Value = -3.180368166004113
$ws.Range("F41").Value = -3.174289534833488
$ws.Range("F42").Value = -3.174390532768089
$ws.Range("F45").Value = -3.191074258706326
$ws.Range("F47").Value = -3.166768976167509
$ws.Range("F48").Value = -3.165408027802204
$ws.Range("F51").Value = -3.186587497102925
$ws.Range("F52").Value = -3.151688875491074
$ws.Range("F54").Value = -3.12364570906453
$ws.Range("F55").Value = -3.117350653507664
$ws.Range("F56").Value = -3.134919428746289
$ws.Range("F57").Value = -3.125833979329375
$ws.Range("F58").Value = -3.088905563172096
$ws.Range("F59").Value = -3.112850189989865
$ws.Range("F60").Value = -3.078894523055742
$ws.Range("F61").Value = -3.06947658972893
$ws.Range("F63").Value = -3.047221488057893
$ws.Range("F64").Value = -3.023985703131832
$ws.Range("F66").Value = -2.989220320896396
$ws.Range("F70").Value = -2.911885431347734
$ws.Range("F71").Value = -2.891289911734929
$ws.Range("F72").Value = -2.871956189347761
$ws.Range("F73").Value = -2.851360669734956
$ws.Range("F81").Value = -3.304686623627517
$ws.Range("F82").Value = -3.343205252912985
$ws.Range("F83").Value = -3.353260808984121
$ws.Range("F84").Value = -3.343529085216978
$ws.Range("F86").Value = -3.285016407825924
$ws.Range("F87").Value = -3.238494550791503
$ws.Range("F92").Value = -3.401811308633782
$ws.Range("F93").Value = -3.427818952172407
$ws.Range("F94").Value = -3.441104537219052
$ws.Range("F95").Value = -3.444931022837255
$ws.Range("F99").Value = -3.271136666724325
$ws.Range("F100").Value = -3.184499247475189
$ws.Range("F103").Value = -3.481788631750945
$ws.Range("F105").Value = -3.509913397379022
$ws.Range("F106").Value = -3.506165876203294
$ws.Range("F107").Value = -3.487689575786036
$ws.Range("F108").Value = -3.451363200358307
$ws.Range("F109").Value = -3.401076388893235
$ws.Range("F110").Value = -3.33401051476334
$ws.Range("F114").Value = -3.525543531140494
$ws.Range("F115").Value = -3.546712523378747
$ws.Range("F116").Value = -3.557952951608451
$ws.Range("F117").Value = -3.547533710650981
$ws.Range("F120").Value = -3.439568878201723
$ws.Range("F122").Value = -3.296831039173977
$ws.Range("F126").Value = -3.584217205658122
$ws.Range("F127").Value = -3.591915395026394
$ws.Range("F128").Value = -3.576977238362117
$ws.Range("F129").Value = -3.553636039509883
$ws.Range("F130").Value = -3.514914831617144
$ws.Range("F133").Value = -3.332441544119829
$ws.Range("F141").Value = -3.537814281497965
$ws.Range("F142").Value = -3.490751663362672
$ws.Range("F143").Value = -3.430231245040901
$ws.Range("F144").Value = -3.35748994001364
$ws.Range("F147").Value = -3.618990597228915
$ws.Range("F148").Value = -3.632401783267247
$ws.Range("F149").Value = -3.630090901319487
$ws.Range("F150").Value = -3.615274694076017
$ws.Range("F151").Value = -3.587027152577486
$ws.Range("F152").Value = -3.553237224535279
$ws.Range("F153").Value = -3.507175842384095
$ws.Range("F155").Value = -3.376656471511474
$ws.Range("F158").Value = -3.6368742688396
$ws.Range("F159").Value = -3.644844249715594
$ws.Range("F160").Value = -3.642525440374336
$ws.Range("F162").Value = -3.599393750084506
$ws.Range("F164").Value = -3.51658666395216
$ws.Range("F166").Value = -3.393282041269236
$ws.Range("F169").Value = -3.652064071701041
$ws.Range("F170").Value = -3.659948418025173
$ws.Range("F171").Value = -3.65456581202152
$ws.Range("F172").Value = -3.63649018031301
$ws.Range("F174").Value = -3.574442766040337
$ws.Range("F175").Value = -3.527947238343799
$ws.Range("F176").Value = -3.472605632180272
$ws.Range("F177").Value = -3.406990570410211
$ws.Range("F180").Value = -3.664234467281588
$ws.Range("F181").Value = -3.668742462440036
$ws.Range("F182").Value = -3.661643363573133
$ws.Range("F184").Value = -3.619419346827894
$ws.Range("F185").Value = -3.583620950180054
$ws.Range("F186").Value = -3.53718309055535
$ws.Range("F187").Value = -3.483605912562465
$ws.Range("F188").Value = -3.418059831777374
$ws.Range("F191").Value = -3.675089869388179
$ws.Range("F193").Value = -3.672395052831875
$ws.Range("F194").Value = -3.653496590425868
$ws.Range("F195").Value = -3.626938738338755
$ws.Range("F196").Value = -3.590978534361095
$ws.Range("F198").Value = -3.49221230342455
$ws.Range("F199").Value = -3.427778729008953
$ws.Range("F202").Value = -3.683404219317499
$ws.Range("F203").Value = -3.685935508373037
$ws.Range("F206").Value = -3.633712069395895
$ws.Range("F210").Value = -3.435809690027914
